$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: score update ---
$ws.Cells.Item(6,6).Value = 4

# --- Rows 68-74: score updates only ---
$ws.Cells.Item(68,7).Value = 1

$ws.Cells.Item(69,6).Value = 1
$ws.Cells.Item(69,7).Value = 1

$ws.Cells.Item(70,6).Value = 2
$ws.Cells.Item(70,7).Value = 1

$ws.Cells.Item(71,7).Value = 2

$ws.Cells.Item(72,7).Value = 2

$ws.Cells.Item(73,7).Value = 1

$ws.Cells.Item(74,7).Value = 2

# --- Rows 75-89: content restructure (old row 75 "Modais de Transporte" flashcard
# is replaced by the content that used to live in row 76, everything below shifts
# up by one row, and new flashcards are appended at rows 83-89) ---

# Row 75
$ws.Cells.Item(75,1).NumberFormat = "@"
$ws.Cells.Item(75,1).Value = '75'
$ws.Cells.Item(75,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(75,3).Value = 'Logística'
$ws.Cells.Item(75,4).Value = 'Classificação dos Modais:'
$ws.Cells.Item(75,5).Value = '<ul> <li><b>Velocidade de Transporte:</b> <ol> <li>Aeroviário</li> <li>Rodoviário</li> <li>Ferroviário</li> <li>Aquaviário</li> <li>Dutoviário</li> </ol></li> <li><b>Disponibilidade:</b> <ol> <li>Rodoviário</li> <li>Ferroviário</li> <li>Aeroviário</li> <li>Aquaviário</li> <li>Dutoviário</li> </ol></li> <li><b>Confiabilidade: </b><ol> <li>Dutoviário</li> <li>Rodoviário</li> <li>Ferroviário</li> <li>Aquaviário</li> <li>Aeroviário</li> </ol></li> <li><b>Capacidade de Carga:</b> <ol> <li>Aquaviário</li> <li>Ferroviário</li> <li>Rodoviário</li> <li>Aeroviário</li> <li>Dutoviário</li> </ol></li> <li><b>Frequência:</b> <ol> <li>Dutoviário</li> <li>Rodoviário</li> <li>Aeroviário</li> <li>Ferroviário</li> <li>Aquaviário</li> </ol></li> </ul>'
$ws.Cells.Item(75,6).Value = 0
$ws.Cells.Item(75,7).Value = 3

# Row 76
$ws.Cells.Item(76,1).Value = 76
$ws.Cells.Item(76,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(76,3).Value = 'Logística'
$ws.Cells.Item(76,4).Value = '<b>Distribuição</b>
<i>Características</i>'
$ws.Cells.Item(76,5).Value = '<ul>
	<li><b>Caonceito:</b> <ul> <li>conjunto de ações voltadas à gestão de materiais, iniciando com a saída do produto do processo produtivo e terminando com a entrega no ponto final de consumo</li> </ul></li>
	<li><b>fatores mais importantes ligados à distribuição</b> <ul> <li>Conferência de cargas;</li> <li>Gestão do frete;</li> <li>Gestão do transporte;</li> <li>Análise e desempenho de indicadores;</li> <li>Gestão de Rotas ou Roteirização.</li> </ul></li>
</ul>






















'
$ws.Cells.Item(76,6).Value = 0
$ws.Cells.Item(76,7).Value = 3

# Row 77
$ws.Cells.Item(77,1).Value = 77
$ws.Cells.Item(77,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(77,3).Value = 'Logística'
$ws.Cells.Item(77,4).Value = '<b>Armazenamento</b>
<i>Características</i>'
$ws.Cells.Item(77,5).Value = '<ul>
	<li>atividades que compreende a armazenagem: <ul> <li>receber</li> <li>carregar</li> <li>descarregar</li> <li>conservar</li> </ul></li>
	<li>quatro pontos principais para que uma empresa decida destinar uma parte de sua área útil à armazenagem, Ballou (1993): <ul> <li>reduzir custos de transporte e produção</li> <li>coordenação de suprimento e demanda</li> <li>auxílio ao processo de produção</li> <li><u>auxílio ao processo de marketing.</u></li> </ul></li>
	<li>funções da armazenagem, Ballou (1993): <ul> <li>Abrigo de produtos</li> <li>Consolidação</li> <li>Transferência e Transbordo</li> <li>Agrupamento</li> </ul></li>
	<li><b>codificação:</b> <ul> <li>catalogar, simplificar, especificar, normatizar e padronizar todo o estoque</li> <li>11 dígitos: <ol> <li>XX - Grupo</li> <li>XX - Classe</li> <li>XXXXXX - Código de identificação</li> <li>X - Dígito de Controle</li> </ol></li> </ul></li>
	<li><b>embalagens:</b> <ul> <li>vantagens: <ul> <li>proteção ao produto <ul> <li>manuseio</li> <li>transporte</li> <li>armazenagem</li> </ul></li> </ul></li> </ul></li>
	<li>ações pelos quais passam os materiais armazenados: <ul> <li>Especificação</li> <li>Simplificação</li> <li>Codificação</li> <li>Padronização</li> <li>Catalogação</li> <li>Normalização</li> </ul></li>
	<li><b>sistemas de armazenamento:</b> <ul> <li>Sistema WMS</li> <li>Racks</li> <li>Mezanino</li> <li>Sistema de carrossel</li> <li>Porta-paletes</li> <li>Flow Rack</li> </ul></li>
</ul>'
$ws.Cells.Item(77,6).Value = 0
$ws.Cells.Item(77,7).Value = 3

# Row 78
$ws.Cells.Item(78,1).Value = 78
$ws.Cells.Item(78,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(78,3).Value = 'Logística'
$ws.Cells.Item(78,4).Value = 'Sistemas Logísticos'
$ws.Cells.Item(78,5).Value = '<ul>
	<li>TMS</li>
	<li>WMS</li>
	<li>Sistema de monitoramento de cargas</li>
	<li>Sistemas de roteirização</li>
	<li>Sistemas de gestão de frotas.</li>
</ul>'
$ws.Cells.Item(78,6).Value = 0
$ws.Cells.Item(78,7).Value = 3

# Row 79
$ws.Cells.Item(79,1).Value = 79
$ws.Cells.Item(79,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(79,3).Value = 'Logística'
$ws.Cells.Item(79,4).Value = '<b>Logística Reversa</b>
<i>Características</i>'
$ws.Cells.Item(79,5).Value = '<ul>
	<li>responsabilidade sobre os resíduos produzidos em decorrência do consumo de bens</li>
	<li>devolução, reciclagem e adequada destinação de produtos pós-venda e pós consumo.</li>
	<li>etapas: <ol> <li>Devolução da embalagem ou resíduo para o comerciante;</li> <li>O comerciante devolve para ao fabricante; e</li> <li>O fabricante destina para reuso, reciclagem ou descarte adequado.</li> </ol></li>
	<li>leis que devem ser cumpridas</li>
	<li>preocupação com a lucratividade e sustentabilidade desse processo</li>
	<li>reversa: <ul> <li>transporte dos produtos nas mãos dos clientes de volta para a empresa</li> </ul></li>
	<li>atividades: <ul> <li>aterro sanitário</li> <li>doação</li> <li>processamento das devoluções</li> <li>reciclagem</li> <li>reembalagem</li> <li>remanufatura</li> <li>revenda</li> <li>revitalização</li> <li>recuperação de cargas roubadas ou perdidas</li> </ul></li>
</ul>'
$ws.Cells.Item(79,6).Value = 0
$ws.Cells.Item(79,7).Value = 2

# Row 80
$ws.Cells.Item(80,1).Value = 80
$ws.Cells.Item(80,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(80,3).Value = 'Logística'
$ws.Cells.Item(80,4).Value = '<b>Logística Verde</b>
<i>Conceito</i>'
$ws.Cells.Item(80,5).Value = 'procedimentos de logística que objetivam a preservação do meio ambiente, que incluem desde a embalagem até o modal de transporte utilizado'
$ws.Cells.Item(80,6).Value = 0
$ws.Cells.Item(80,7).Value = 2

# Row 81
$ws.Cells.Item(81,1).Value = 81
$ws.Cells.Item(81,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(81,3).Value = 'Logística'
$ws.Cells.Item(81,4).Value = '<b>Logística de Pós-consumo</b>
<i>Características</i>'
$ws.Cells.Item(81,5).Value = '<ul>
	<li>favorece o retorno dos produtos após serem utilizados pelos clientes, visando: <ul> <li>reciclagem</li> <li>reutilização ou</li> <li>descarte apropriado</li> </ul></li>
	<li><b>motivação:</b> <ul> <li>quantidade de materiais descartados pela sociedade desde o século XX até os dias de hoje <ul> <li>Diminuilção do ciclo de vida dos produtos</li> </ul> </li> </ul></li>
	<li>destinos para um produto após descarte: <ul> <li>local seguro (aterro sanitário)</li> <li>local não seguro</li> <li>Reciclagem</li> </ul></li>
</ul>'
$ws.Cells.Item(81,6).Value = 0
$ws.Cells.Item(81,7).Value = 2

# Row 82
$ws.Cells.Item(82,1).Value = 82
$ws.Cells.Item(82,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(82,3).Value = 'Logística'
$ws.Cells.Item(82,4).Value = '<b>Logística de Pós-venda</b>
<i>Características</i>'
$ws.Cells.Item(82,5).Value = '<b>uma das suas preocupações-chave: </b><ul> <li>criar um canal acessível para clientes retornarem <b>produtos</b>. <ul> <li>defeitos de fabricação ou</li> <li>erros no pedido</li> </ul></li> </ul>'
$ws.Cells.Item(82,6).Value = 0
$ws.Cells.Item(82,7).Value = 1

# Row 83
$ws.Cells.Item(83,1).NumberFormat = "@"
$ws.Cells.Item(83,1).Value = '83'
$ws.Cells.Item(83,2).Value = 'Matemática Básica'
$ws.Cells.Item(83,3).Value = 'Logarítimos'
$ws.Cells.Item(83,4).Value = 'Propriedades dos logaritmos'
$ws.Cells.Item(83,5).Value = '<ul>
	<li>log<sub>a</sub>(b &sdot; c) =log<sub>a</sub>b + log<sub>a</sub>c</li>
	<li>log<sub>a</sub>(b/c) = log<sub>a</sub>b - log<sub>a</sub>c</li>
	<li>log<sub>a</sub>b<sup>c</sup> = clog<sub>a</sub>b</li>
	<li>log<sub>a<sup>x</sup></sub>b = (1/x)log<sub>a</sub>b</li>
	<li>log<sub>b</sub>a = (log<sub>c</sub>a)/(log<sub>c</sub>b)</li>
</ul>'
$ws.Cells.Item(83,6).Value = 1
$ws.Cells.Item(83,7).Value = 1

# Row 84
$ws.Cells.Item(84,1).Value = 84
$ws.Cells.Item(84,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(84,3).Value = 'Logística'
$ws.Cells.Item(84,4).Value = '<b>Modais de Transporte</b>: Aeroviário'
$ws.Cells.Item(84,5).Value = '<ul>
	<li>vantagens:<ul> <li>longas distâncias, independente de acidentes ou formações geográficas</li> <li>mais rápido dentre os modais</li> <li>Menor custo com embalagens</li> </ul></li>
	<li>desvantagens:<ul> <li>volume pequeno de cargas</li> <li>custo mais elevado</li> <li>grande chance de precisar de outro modal para que o produto chegue ao destino</li> </ul></li>
</ul>'
$ws.Cells.Item(84,6).Value = 0
$ws.Cells.Item(84,7).Value = 1

# Row 85
$ws.Cells.Item(85,1).Value = 85
$ws.Cells.Item(85,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(85,3).Value = 'Logística'
$ws.Cells.Item(85,4).Value = '<b>Modais de transporte</b>: Aquaviário'
$ws.Cells.Item(85,5).Value = '<ul>
	<li>Marítimo: mares e oceanos;</li>
	<li>Fluvial: rios;</li>
	<li>Lacustre: lagos e lagoas.</li>
	<li>vantagens: <ul> <li>Maior capacidade de carga</li> <li>grandes distâncias de forma autônoma</li> <li>Baixo custo unitário de carregamento</li> </ul></li>
	<li>desvantagens: <ul> <li>é o mais lento entre os modais</li> <li>Maior suscetibilidade as mudanças da natureza</li> <li>Necessidade de terminais especializados</li> <li>Desembaraço burocrático</li> <li>Alto custo quanto ao seguro das cargas</li> </ul></li>
</ul>'
$ws.Cells.Item(85,6).Value = 0
$ws.Cells.Item(85,7).Value = 1

# Row 86
$ws.Cells.Item(86,1).Value = 86
$ws.Cells.Item(86,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(86,3).Value = 'Logística'
$ws.Cells.Item(86,4).Value = '<b>Modasi de Transporte</b>: Ferroviário'
$ws.Cells.Item(86,5).Value = '<ul>
	<li>vantagens: <ul> <li>Baixo custo</li> <li>Menor risco de acidentes</li> <li>grande capacidade de transporte de cargas</li> </ul></li>
	<li>desvantagens: <ul> <li>grande chance de precisar de outro modal para que o produto chegue ao destino;</li> <li>Baixo investimento governamental</li> <li>Rotas fixas e inflexíveis</li> </ul></li>
</ul>'
$ws.Cells.Item(86,6).Value = 0
$ws.Cells.Item(86,7).Value = 1

# Row 87
$ws.Cells.Item(87,1).Value = 87
$ws.Cells.Item(87,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(87,3).Value = 'Logística'
$ws.Cells.Item(87,4).Value = '<b>Modais de Transporte</b>: Rodoviário'
$ws.Cells.Item(87,5).Value = '<ul>
	<li>realizado por meio de caminhões, carretas, carros e veículos em geral através de rodovias.</li>
	<li>vantagens: <ul> <li>alcance (acessibilidade)</li> <li>Rapidez para contratação</li> <li>Rotas flexíveis</li> <li>Menor burocracia em relação aos outros modais</li> <li>custo de estrutura é menor e há investimento governamental</li> </ul></li>
	<li>desvantagens: <ul> <li>Gastos com pedágio e maior suscetibilidade ao aumento de combustíveis tendem a aumentar o valor do frete</li> <li>capacidade de carga é bem menor</li> <li>autonomia de deslocamento é menor</li> <li>Maior chance de extravio quanto a carga</li> </ul></li>
</ul>'
$ws.Cells.Item(87,6).Value = 0
$ws.Cells.Item(87,7).Value = 1

# Row 88
$ws.Cells.Item(88,1).Value = 88
$ws.Cells.Item(88,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(88,3).Value = 'Logística'
$ws.Cells.Item(88,4).Value = '<b>Modais de Transporte</b>: Dutoviário'
$ws.Cells.Item(88,5).Value = '<ul>
	<li>largamente utilizado na indústria de Petróleo e seus derivados</li>
	<li>vantagens: <ul> <li>Consegue transportar grande volume de carga de forma constante</li> <li>grande confiabilidade no processo de captação e entrega</li> <li>Baixos custos operacionais</li> </ul></li>
	<li>desvantagens: <ul> <li>Custo inicial para implantação do projeto (fixo) altíssimo</li> <li>Burocracia ambiental</li> <li>Reduzida flexibilidade de trajeto</li> </ul></li>
</ul>'
$ws.Cells.Item(88,6).Value = 0
$ws.Cells.Item(88,7).Value = 1

# Row 89
$ws.Cells.Item(89,1).Value = 89
$ws.Cells.Item(89,2).Value = 'Conhecimentos Específicos'
$ws.Cells.Item(89,3).Value = 'Logística'
$ws.Cells.Item(89,4).Value = '<b>Cabotagem</b>
<i>Conceito</i>'
$ws.Cells.Item(89,5).Value = 'transporte marítimo ou hidroviário de cargas e passageiros entre portos ou pontos dentro do mesmo país, mantendo a costa à vista, diferente da navegação de longo curso (internacional)'
$ws.Cells.Item(89,6).Value = 0
$ws.Cells.Item(89,7).Value = 1

